# Auto-generated Excel COM-interop script to update cryptos list
# matching the commit 'Updated cryptos list on Wed May 29 02:05:37 UTC 2024 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "68.592.42"
Set-TextValue $ws.Range("E2") "  +0.06%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.841.80"
Set-TextValue $ws.Range("E3") "  -0.02%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.13%  "

# Row 5
Set-TextValue $ws.Range("D5") "602.49"
Set-TextValue $ws.Range("E5") "  +0.28%  "

# Row 6
Set-TextValue $ws.Range("D6") "168.42"
Set-TextValue $ws.Range("E6") "  +0.45%  "

# Row 7
Set-TextValue $ws.Range("D7") "3.841.93"
Set-TextValue $ws.Range("E7") "  +0.11%  "

# Row 8
Set-TextValue $ws.Range("E8") "  -0.13%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.525"
Set-TextValue $ws.Range("E9") "  -0.57%  "

# Row 10
Set-TextValue $ws.Range("E10") "  +0.82%  "

# Row 11
Set-TextValue $ws.Range("D11") "6.47"
Set-TextValue $ws.Range("E11") "  +2.33%  "

# Row 12
Set-TextValue $ws.Range("E12") "  -0.54%  "

# Row 13
Set-TextValue $ws.Range("D13") "0.0000279"
Set-TextValue $ws.Range("E13") "  +11.57%  "

# Row 14
Set-TextValue $ws.Range("D14") "36.83"
Set-TextValue $ws.Range("E14") "  -1.64%  "

# Row 15
Set-TextValue $ws.Range("D15") "4.485.37"
Set-TextValue $ws.Range("E15") "  -0.55%  "

# Row 16
Set-TextValue $ws.Range("D16") "3.840.28"

# Row 17
Set-TextValue $ws.Range("D17") "68.488.74"
Set-TextValue $ws.Range("E17") "  -0.31%  "

# Row 18
Set-TextValue $ws.Range("D18") "18.38"
Set-TextValue $ws.Range("E18") "  +0.13%  "

# Row 19
Set-TextValue $ws.Range("D19") "7.34"
Set-TextValue $ws.Range("E19") "  -2.26%  "

# Row 20
Set-TextValue $ws.Range("E20") "  -0.46%  "

# Row 21
Set-TextValue $ws.Range("D21") "11.02"
Set-TextValue $ws.Range("E21") "  +2.13%  "

# Row 22
Set-TextValue $ws.Range("D22") "469.92"
Set-TextValue $ws.Range("E22") "  -1.96%  "

# Row 23
Set-TextValue $ws.Range("D23") "0.724"
Set-TextValue $ws.Range("E23") "  -1.14%  "

# Row 24
Set-TextValue $ws.Range("D24") "0.0000162"
Set-TextValue $ws.Range("E24") "  +0.74%  "

# Row 25
Set-TextValue $ws.Range("D25") "83.29"
Set-TextValue $ws.Range("E25") "  -1.68%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.24"
Set-TextValue $ws.Range("E26") "  -0.18%  "

# Row 27
Set-TextValue $ws.Range("D27") "12.05"
Set-TextValue $ws.Range("E27") "  -1.38%  "

# Row 28
Set-TextValue $ws.Range("D28") "10.39"
Set-TextValue $ws.Range("E28") "  +4.16%  "

# Row 29
Set-TextValue $ws.Range("E29") "  +0.07%  "

# Row 30
Set-TextValue $ws.Range("D30") "2.95"
Set-TextValue $ws.Range("E30") "  -0.34%  "

# Row 31
Set-TextValue $ws.Range("D31") "3.991.54"
Set-TextValue $ws.Range("E31") "  -0.24%  "

# Row 32
Set-TextValue $ws.Range("D32") "7.67"
Set-TextValue $ws.Range("E32") "  -0.80%  "

# Row 33
Set-TextValue $ws.Range("D33") "31.21"
Set-TextValue $ws.Range("E33") "  +0.40%  "

# Row 34
Set-TextValue $ws.Range("D34") "2.29"
Set-TextValue $ws.Range("E34") "  -0.91%  "

# Row 35
Set-TextValue $ws.Range("D35") "9.26"
Set-TextValue $ws.Range("E35") "  -1.94%  "

# Row 36
Set-TextValue $ws.Range("D36") "3.808.90"
Set-TextValue $ws.Range("E36") "  -0.33%  "

# Row 37
Set-TextValue $ws.Range("E37") "  -0.94%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.72"
Set-TextValue $ws.Range("E38") "  +12.87%  "

# Row 39
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D39") "1.02"
Set-TextValue $ws.Range("E39") "  +0.19%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D40") "0.139"
Set-TextValue $ws.Range("E40") "  +0.18%  "

# Row 41
Set-TextValue $ws.Range("D41") "5.90"
Set-TextValue $ws.Range("E41") "  -1.35%  "

# Row 42
Set-TextValue $ws.Range("E42") "  +0.05%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.316"
Set-TextValue $ws.Range("E43") "  -0.08%  "

# Row 44
Set-TextValue $ws.Range("D44") "2.00"
Set-TextValue $ws.Range("E44") "  -1.73%  "

# Row 45
Set-TextValue $ws.Range("D45") "420.89"
Set-TextValue $ws.Range("E45") "  -0.73%  "

# Row 46
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D46") "1.00"
Set-TextValue $ws.Range("E46") "  -0.01%  "

# Row 47
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D47") "8.66"
Set-TextValue $ws.Range("E47") "  +1.02%  "

# Row 48
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue $ws.Range("D48") "0.000291"
Set-TextValue $ws.Range("E48") "  +8.38%  "

# Row 49
Set-TextValue $ws.Range("D49") "46.92"
Set-TextValue $ws.Range("E49") "  -1.59%  "

# Row 50
Set-TextValue $ws.Range("D50") "141.58"
Set-TextValue $ws.Range("E50") "  +0.65%  "

# Row 51
Set-TextValue $ws.Range("D51") "26.00"
Set-TextValue $ws.Range("E51") "  +4.16%  "
